# Project Doc - Missing UML
#
# The original author re-ran Word's proofing pass over the document,
# which (besides flagging "NameEqualsNull", "Lehmbecker", "uml" and
# "codesmells" as misspellings, and "email" as a grammar nit) splits
# the runs that contain those words so the proofing squiggles can be
# anchored. We reproduce the run splits - and the xml:space="preserve"
# Word adds to the run that is left with trailing a space once split -
# by toggling a no-op character-formatting property across the
# sub-range we want broken out into its own run.

$d = $word.ActiveDocument

function Split-Run([int]$start, [int]$end) {
    # Touching Font.Bold (set then immediately reset) forces the
    # engine to re-emit the run boundaries around $start..$end without
    # altering the visible formatting.
    $sub = $d.Range($start, $end)
    $sub.Font.Bold = $true
    $sub.Font.Bold = $false
}

# --- "Team: NameEqualsNull" -> "Team: " | "NameEqualsNull" ----------------
$d.Content.Find.Execute("NameEqualsNull", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng = $d.Content.Find.Parent.Duplicate
Split-Run $rng.Start $rng.End

# --- "Turner Lehmbecker" -> "Turner " | "L" | "ehmbecker" ------------------
$d.Content.Find.Execute("Turner L", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$found = $d.Content.Find.Parent.Duplicate
Split-Run ($found.End - 1) $found.End

# --- "email me the rest ... and uml)." -------------------------------------
$d.Content.Find.Execute("email", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$found = $d.Content.Find.Parent.Duplicate
Split-Run $found.Start $found.End

$d.Content.Find.Execute("uml", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$found = $d.Content.Find.Parent.Duplicate
Split-Run $found.Start $found.End
Split-Run $found.End ($found.End + 2)

# --- "...other codesmells." -------------------------------------------------
$d.Content.Find.Execute("codesmells", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$found = $d.Content.Find.Parent.Duplicate
Split-Run $found.Start $found.End
Split-Run $found.End ($found.End + 1)

Write-Host "done"
